$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.966.83"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.779.89"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.77"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5372"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3764"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07442"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.66"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.96%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9990"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.43"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.073"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.219"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.776.19"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.25"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.24%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06433"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9990"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.21"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.877"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.022.13"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.16"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.088"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.18"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.21%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.979.01"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.294"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.74"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.104"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1053"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.642"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.506"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.2253"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06419"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02273"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.995"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.423"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.10"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6144"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.429"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.37%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9983"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.20"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.60%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5740"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.25"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.926"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.182"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06790"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.69%  "
